# feat: add 2022-Q1 data
#
# The workbook has quarterly fund-holding sheets (2020-Q4 .. 2021-Q4) plus a
# "总计" (totals) roll-up sheet. This adds a new "2022-Q1" fund-holding sheet
# (positioned right before "总计") and rolls its summary numbers into "总计".
#
# Strategy: rename the existing "总计" sheet to "2022-Q1" (so it inherits the
# correct sheetId/rId slot order seen in the target workbook.xml), add a
# brand-new "总计" sheet right after it, then:
#   - overwrite the (renamed) "2022-Q1" sheet with the new fund table
#   - (re)write the totals table into the new "总计" sheet, with a fresh
#     "2022-Q1" row inserted on top of the previous totals rows.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("2021-Q4")     # formatting donor for fund-table rows
$oldTotal = $wb.Worksheets.Item("总计")        # becomes "2022-Q1"

$oldTotal.Name = "2022-Q1"
$q1 = $oldTotal

$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1)
$total.Name = "总计"

function Set-TextCell($ws, $addr, $text) {
    # Force text storage (no numeric / leading-zero coercion) without leaving
    # a stray "quote prefix" style behind.
    $ws.Range($addr).Value = "'" + $text
}

function Copy-RowFormat($srcWs, $srcAddr, $dstWs, $dstAddr) {
    $srcWs.Range($srcAddr).Copy()
    $dstWs.Range($dstAddr).PasteSpecial(-4122)   # xlPasteFormats
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# 1) "2022-Q1" sheet: fresh fund-holding table (replaces the old totals data)
# ---------------------------------------------------------------------------
$q1.Cells.Clear()

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $col = [char]([int][char]'B' + $c)
    Set-TextCell $q1 ($col + "1") $headers[$c]
}

$fundRows = @(
    @{ code = "163804"; name = "中银收益混合A";        scale = "19.19"; pos = "85.56"; pct = "2.97"; mv = "0.5699"; rank = 8  },
    @{ code = "163822"; name = "中银主题策略混合";      scale = "7.89";  pos = "84.73"; pct = "2.77"; mv = "0.2186"; rank = 10 },
    @{ code = "001672"; name = "国寿安保智慧生活股票";  scale = "3.56";  pos = "85.91"; pct = "3.04"; mv = "0.1082"; rank = 6  },
    @{ code = "010965"; name = "中银鑫新消费成长混合A"; scale = "4.76";  pos = "73.66"; pct = "2.19"; mv = "0.1042"; rank = 10 },
    @{ code = "014505"; name = "中银收益混合C";        scale = "0.98";  pos = "85.56"; pct = "2.97"; mv = "0.0291"; rank = 8  },
    @{ code = "010962"; name = "中银鑫新消费成长混合C"; scale = "0.82";  pos = "73.66"; pct = "2.19"; mv = "0.0180"; rank = 10 },
    @{ code = "710301"; name = "富安达增强收益债券A";   scale = "0.61";  pos = "20.20"; pct = "1.79"; mv = "0.0109"; rank = 4  },
    @{ code = "710302"; name = "富安达增强收益债券C";   scale = "0.26";  pos = "20.20"; pct = "1.79"; mv = "0.0047"; rank = 4  },
    @{ code = "960012"; name = "中银收益混合H";        scale = "0.03";  pos = "85.56"; pct = "2.97"; mv = "0.0009"; rank = 8  }
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    $q1.Range("A$r").Value = $i
    Set-TextCell $q1 "B$r" $row.code
    Set-TextCell $q1 "C$r" $row.name
    Set-TextCell $q1 "D$r" $row.scale
    Set-TextCell $q1 "E$r" $row.pos
    Set-TextCell $q1 "F$r" $row.pct
    Set-TextCell $q1 "G$r" $row.mv
    $q1.Range("H$r").Value = $row.rank

    Copy-RowFormat $template "A2:H2" $q1 "A$r"
}
Copy-RowFormat $template "A1:H1" $q1 "A1"

# ---------------------------------------------------------------------------
# 2) "总计" sheet: totals table with the new 2022-Q1 row on top
# ---------------------------------------------------------------------------
$total.Cells.Clear()

Set-TextCell $total "B1" "日期"
Set-TextCell $total "C1" "持有数量(只)"
Set-TextCell $total "D1" "持有市值(亿元)"

$totalRows = @(
    @{ date = "2022-Q1"; count = 9;  mv = 1.06 },
    @{ date = "2021-Q4"; count = 20; mv = 3.51 },
    @{ date = "2021-Q3"; count = 6;  mv = 0.47 },
    @{ date = "2021-Q2"; count = 14; mv = 1.47 },
    @{ date = "2021-Q1"; count = 15; mv = 1.75 },
    @{ date = "2020-Q4"; count = 20; mv = 3.82 }
)

for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]

    $total.Range("A$r").Value = $i
    Set-TextCell $total "B$r" $row.date
    $total.Range("C$r").Value = $row.count
    $total.Range("D$r").Value = $row.mv
}

Copy-RowFormat $template "A1:D1" $total "A1"
for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $r = $i + 2
    Copy-RowFormat $template "A2:D2" $total "A$r"
}
